$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header D1: "WebExcel" -> "ExDataExcel"
$ws.Cells.Item(1, 4).Value = "ExDataExcel"

# Update data rows D2:D4: "...WebData.xlsx" -> "...ExcelData.xlsx"
$ws.Cells.Item(2, 4).Value = "TC01_CDSValidation_by_ParticipantID - 1_ExcelData.xlsx"
$ws.Cells.Item(3, 4).Value = "TC01_CDSValidation_by_ParticipantID - 1_ExcelData.xlsx"
$ws.Cells.Item(4, 4).Value = "TC01_CDSValidation_by_ParticipantID - 1_ExcelData.xlsx"

# Update selection to D2
$ws.Range("D2").Select()
